# Update cryptos list with latest prices/volumes/links
# Applies the row-level changes captured in the target OOXML diff for cryptos.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='30.704.20'; E='  -1.08%  '},
    @{Row=3; D='1.924.53'; E='  -1.63%  '},
    @{Row=4; D='1.001'; E='  +0.17%  '; DNumeric=$true},
    @{Row=5; D='241.70'; E='  -1.44%  '; DNumeric=$true},
    @{Row=6; D='1.001'; E='  +0.12%  '; DNumeric=$true},
    @{Row=7; D='0.4846'; E='  -0.53%  '; DNumeric=$true},
    @{Row=8; D='0.2922'; E='  -1.22%  '; DNumeric=$true},
    @{Row=9; D='0.06802'; E='  -0.40%  '; DNumeric=$true},
    @{Row=10; E='  -0.61%  '},
    @{Row=11; D='105.92'; E='  -1.50%  '; DNumeric=$true},
    @{Row=12; D='1.935.17'; E='  -1.01%  '},
    @{Row=13; D='0.07757'; E='  -0.58%  '; DNumeric=$true},
    @{Row=14; D='5.307'; E='  -2.76%  '; DNumeric=$true},
    @{Row=15; D='0.6943'; E='  -1.55%  '; DNumeric=$true},
    @{Row=16; D='273.12'; E='  -4.23%  '; DNumeric=$true},
    @{Row=17; D='30.719.73'; E='  -1.09%  '},
    @{Row=18; D='0.000007648'; E='  -0.64%  '; DNumeric=$true},
    @{Row=19; E='  +0.07%  '},
    @{Row=20; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.557'; E='  +0.72%  '; DNumeric=$true},
    @{Row=21; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='12.91'; E='  -2.21%  '; DNumeric=$true},
    @{Row=22; B='BitDAO'; C='https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'; D='0.4609'; E='  -6.46%  '; DNumeric=$true},
    @{Row=23; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.001'; E='  +0.07%  '; DNumeric=$true},
    @{Row=24; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='6.425'; E='  -1.21%  '; DNumeric=$true},
    @{Row=25; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='9.807'; E='  +0.03%  '; DNumeric=$true},
    @{Row=26; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='164.57'; E='  -3.50%  '; DNumeric=$true},
    @{Row=27; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='19.37'; E='  -3.14%  '; DNumeric=$true},
    @{Row=28; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='2.140'; E='  -2.81%  '; DNumeric=$true},
    @{Row=29; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.1034'; E='  -2.22%  '; DNumeric=$true},
    @{Row=30; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.379'; E='  -2.02%  '; DNumeric=$true},
    @{Row=31; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='4.555'; E='  -1.41%  '; DNumeric=$true},
    @{Row=32; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='1.542'; E='  -2.74%  '; DNumeric=$true},
    @{Row=33; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='4.330'; E='  -3.02%  '; DNumeric=$true},
    @{Row=34; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.04846'; E='  -1.62%  '; DNumeric=$true},
    @{Row=35; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.7530'; E='  -1.42%  '; DNumeric=$true},
    @{Row=36; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.135'; E='  -3.13%  '; DNumeric=$true},
    @{Row=37; B='Frax'; C='https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D='1.001'; E='  +0.05%  '; DNumeric=$true},
    @{Row=38; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.718'; E='  -0.40%  '; DNumeric=$true},
    @{Row=39; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.01981'; E='  -1.46%  '; DNumeric=$true},
    @{Row=40; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.645'; E='  -2.18%  '; DNumeric=$true},
    @{Row=41; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='6.446'; E='  -1.22%  '; DNumeric=$true},
    @{Row=42; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='76.74'; E='  +1.64%  '; DNumeric=$true},
    @{Row=43; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.044'; E='  -2.68%  '; DNumeric=$true},
    @{Row=44; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.8778'; E='  -1.84%  '; DNumeric=$true},
    @{Row=45; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.4399'; E='  -1.73%  '; DNumeric=$true},
    @{Row=46; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='107.01'; E='  -2.11%  '; DNumeric=$true},
    @{Row=47; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='7.840'; E='  -4.17%  '; DNumeric=$true},
    @{Row=48; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='1.000'; E='  +0.02%  '; DNumeric=$true},
    @{Row=49; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='973.66'; E='  -2.55%  '; DNumeric=$true},
    @{Row=50; D='0.1228'; E='  -2.31%  '; DNumeric=$true},
    @{Row=51; B='Elrond'; C='https://coinranking.com/coin/omwkOTglq+elrond-egld'; D='35.95'; E='  +0.55%  '; DNumeric=$true}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        $dCell = $ws.Cells.Item($r, 4)
        if ($u.ContainsKey("DNumeric")) {
            # Value looks like a plain number (e.g. "1.001"); force text
            # storage so Excel does not collapse it to a Double and lose
            # the exact original formatting (trailing zeros, etc.)
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}

Write-Host "Applied $($updates.Count) row updates"
